# Serbia Prva Liga - Atualização de bases das ligas, do dia: 10-06-2024 às 07:08
#
# The edit consists of re-pairing certain rows of match data (columns B..AD)
# while leaving the row-index column A untouched. This is implemented as a
# series of row-content swaps (and two 3-way rotations) using the data range
# B<row>:AD<row> for each affected row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $r1, $r2) {
    $rng1 = $ws.Range("B$r1`:AD$r1")
    $rng2 = $ws.Range("B$r2`:AD$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

# Rotate content forward through a cycle of rows such that:
#   new(rows[0]) = old(rows[-1])
#   new(rows[i]) = old(rows[i-1])  for i > 0
function Rotate-Rows($ws, [int[]]$rows) {
    $ranges = @()
    $values = @()
    foreach ($r in $rows) {
        $rng = $ws.Range("B$r`:AD$r")
        $ranges += $rng
        $values += ,$rng.Value2
    }
    $n = $ranges.Count
    for ($i = 0; $i -lt $n; $i++) {
        $prev = ($i - 1 + $n) % $n
        $ranges[$i].Value2 = $values[$prev]
    }
}

# Simple pairwise swaps
Swap-Rows $ws 19 20
Swap-Rows $ws 140 141
Swap-Rows $ws 153 154
Swap-Rows $ws 158 159
Swap-Rows $ws 197 198
Swap-Rows $ws 212 213
Swap-Rows $ws 219 221
Swap-Rows $ws 248 250
Swap-Rows $ws 253 254

# 3-way rotations
Rotate-Rows $ws @(155, 156, 157)
Rotate-Rows $ws @(171, 172, 173)
